# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (per-fund holding detail) right before
#   the existing "总计" (totals) sheet.
# - Add a new summary row for "2022-Q1" at the top of the "总计" sheet's
#   data (existing rows shift down).

$wb = $excel.ActiveWorkbook

# Helper: write a value as genuine TEXT (not auto-converted to a number),
# then drop the cell back to the workbook's default ("Normal") style so we
# don't leave a stray custom number format behind on a plain data cell.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Helper: apply the bold / centered / thin-bordered look used for header
# row + row-index column throughout this workbook.
function Set-HeaderLook($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# ------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet just before "总计".
# ------------------------------------------------------------------------
# The "总计" sheet is currently the last (4th) tab; insert the new sheet
# immediately before it so the tab order becomes:
#   2020-Q4, 2021-Q3, 2021-Q4, 2022-Q1, 总计
$totalSheetBeforeInsert = $wb.Worksheets.Item(4)
$newSheet = $wb.Worksheets.Add($totalSheetBeforeInsert)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet handles returned by Item(...) track tab *position*, not
# sheet identity - after the insert above, the old $totalSheetBeforeInsert
# handle now resolves to the newly-added sheet instead. Re-fetch "总计" by
# name now that it has moved to position 5.
$totalSheet = $wb.Worksheets.Item("总计")

# ---- "2022-Q1" header row ----------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
Set-HeaderLook $newSheet.Range("B1:H1")

# ---- "2022-Q1" data rows ------------------------------------------------
# Row-index column (A) and "仓位排名" column (H) are numbers; every other
# data column is stored as text, matching the other quarterly sheets.
Set-TextValue $newSheet.Range("B2") "010404"
Set-TextValue $newSheet.Range("C2") "博道盛利6个月持有期混合"
Set-TextValue $newSheet.Range("D2") "1.29"
Set-TextValue $newSheet.Range("E2") "34.13"
Set-TextValue $newSheet.Range("F2") "0.59"
Set-TextValue $newSheet.Range("G2") "0.0076"
$newSheet.Range("H2").Value = 5

Set-TextValue $newSheet.Range("B3") "164811"
Set-TextValue $newSheet.Range("C3") "工银瑞信中证京津冀协同发展主题指数（LOF）A"
Set-TextValue $newSheet.Range("D3") "0.23"
Set-TextValue $newSheet.Range("E3") "94.28"
Set-TextValue $newSheet.Range("F3") "2.94"
Set-TextValue $newSheet.Range("G3") "0.0068"
$newSheet.Range("H3").Value = 9

Set-TextValue $newSheet.Range("B4") "512780"
Set-TextValue $newSheet.Range("C4") "广发中证京津冀协同发展主题ETF"
Set-TextValue $newSheet.Range("D4") "0.13"
Set-TextValue $newSheet.Range("E4") "98.52"
Set-TextValue $newSheet.Range("F4") "3.14"
Set-TextValue $newSheet.Range("G4") "0.0041"
$newSheet.Range("H4").Value = 7

Set-TextValue $newSheet.Range("B5") "164825"
Set-TextValue $newSheet.Range("C5") "工银瑞信中证京津冀协同发展主题指数（LOF）C"
Set-TextValue $newSheet.Range("D5") "0.06"
Set-TextValue $newSheet.Range("E5") "94.28"
Set-TextValue $newSheet.Range("F5") "2.94"
Set-TextValue $newSheet.Range("G5") "0.0018"
$newSheet.Range("H5").Value = 9

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
Set-HeaderLook $newSheet.Range("A2:A5")

# ------------------------------------------------------------------------
# 2) Update "总计": insert a new top data row for 2022-Q1, existing rows
#    (2021-Q4 / 2021-Q3 / 2020-Q4) shift down by one.
# ------------------------------------------------------------------------
$totalSheet.Rows(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.02
$totalSheet.Range("B2:D2").Style = "Normal"
Set-HeaderLook $totalSheet.Range("A2")

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

Write-Host "2022-Q1 sheet added and total sheet updated"
